# Apply content edits to slide 2 ("Experiment 3" flow diagram):
#   - "Self" -> "Mirror" (two occurrences)
#   - "Antiself" -> "Inverse" (two occurrences)
#   - "Analysis: n = 89" -> "Analysis: n = 97" (first occurrence only)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$s.Shapes.Item(31).TextFrame.TextRange.Text = "Mirror"
$s.Shapes.Item(32).TextFrame.TextRange.Text = "Inverse"
$s.Shapes.Item(35).TextFrame.TextRange.Text = "Mirror"
$s.Shapes.Item(36).TextFrame.TextRange.Text = "Inverse"
$s.Shapes.Item(46).TextFrame.TextRange.Text = "Analysis: n = 97"
